$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values stay as text (they include formats like
# thousand-dot-separated numbers, trailing zeros, scientific-looking decimals, etc.)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '72.030.22'
$ws.Range("D3").Value = '3.626.26'
$ws.Range("E3").Value = '  +7.01%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '597.73'
$ws.Range("E5").Value = '  +1.67%  '
$ws.Range("D6").Value = '182.69'
$ws.Range("E6").Value = '  +1.58%  '
$ws.Range("D7").Value = '3.612.67'
$ws.Range("E7").Value = '  +6.86%  '
$ws.Range("E8").Value = '  +1.98%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("E10").Value = '  +6.03%  '
$ws.Range("E11").Value = '  +3.10%  '
$ws.Range("D12").Value = '50.57'
$ws.Range("E12").Value = '  +4.34%  '
$ws.Range("D13").Value = '0.0000290'
$ws.Range("E13").Value = '  +2.86%  '
$ws.Range("D14").Value = '706.20'
$ws.Range("E14").Value = '  +4.01%  '
$ws.Range("D15").Value = '4.202.36'
$ws.Range("E15").Value = '  +6.99%  '
$ws.Range("E16").Value = '  +3.75%  '
$ws.Range("D17").Value = '72.103.43'
$ws.Range("E17").Value = '  +3.90%  '
$ws.Range("D18").Value = '3.625.41'
$ws.Range("E18").Value = '  +7.41%  '
$ws.Range("D19").Value = '0.123'
$ws.Range("E19").Value = '  +1.76%  '
$ws.Range("D20").Value = '18.56'
$ws.Range("E20").Value = '  +4.96%  '
$ws.Range("E21").Value = '  +4.48%  '
$ws.Range("D22").Value = '0.935'
$ws.Range("E22").Value = '  +3.42%  '
$ws.Range("D23").Value = '5.83'
$ws.Range("E23").Value = '  +7.46%  '
$ws.Range("D24").Value = '17.89'
$ws.Range("E24").Value = '  +4.41%  '
$ws.Range("D25").Value = '105.43'
$ws.Range("E25").Value = '  +2.26%  '
$ws.Range("E26").Value = '  +2.91%  '
$ws.Range("E27").Value = '  +4.69%  '
$ws.Range("D28").Value = '10.04'
$ws.Range("E28").Value = '  +4.32%  '
$ws.Range("D29").Value = '35.72'
$ws.Range("E29").Value = '  +5.14%  '
$ws.Range("E30").Value = '  +4.36%  '
$ws.Range("D31").Value = '7.42'
$ws.Range("E31").Value = '  +7.11%  '
$ws.Range("D32").Value = '4.19'
$ws.Range("E32").Value = '  +19.07%  '
$ws.Range("D33").Value = '592.55'
$ws.Range("E33").Value = '  +6.10%  '
$ws.Range("D34").Value = '11.39'
$ws.Range("E34").Value = '  +2.33%  '
$ws.Range("E35").Value = '  +1.52%  '
$ws.Range("D36").Value = '59.83'
$ws.Range("E36").Value = '  +2.04%  '
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '0.145'
$ws.Range("E38").Value = '  +4.33%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '3.653.28'
$ws.Range("E39").Value = '  -0.44%  '
$ws.Range("D40").Value = '36.07'
$ws.Range("E40").Value = '  +1.02%  '
$ws.Range("D41").Value = '0.0₃0776'
$ws.Range("E41").Value = '  +10.54%  '
$ws.Range("D42").Value = '3.44'
$ws.Range("E42").Value = '  +5.05%  '
$ws.Range("E43").Value = '  +4.49%  '
$ws.Range("D44").Value = '0.0453'
$ws.Range("E44").Value = '  +7.16%  '
$ws.Range("E45").Value = '  +2.26%  '
$ws.Range("E46").Value = '  +2.37%  '
$ws.Range("D47").Value = '2.80'
$ws.Range("E47").Value = '  +4.52%  '
$ws.Range("D48").Value = '1.48'
$ws.Range("E48").Value = '  +4.74%  '
$ws.Range("E49").Value = '  +2.11%  '
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("D51").Value = '133.42'
$ws.Range("E51").Value = '  +0.07%  '
